$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename distribution title, set Payments Paid to Yes, clear
# Distribution Basis / From Currency / To Currency / Exchange Rate / As Of
$ws.Range("C2").Value = "Distribution 1"
$ws.Range("I2").Value = "Yes"
$ws.Range("K2").Value = $null
$ws.Range("L2").Value = $null
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null
$ws.Range("O2").Value = $null

# Row 3: rename distribution title, set Payments Paid to Yes, clear
# Distribution Basis / From Currency / To Currency / Exchange Rate / As Of
$ws.Range("C3").Value = "Distribution 2"
$ws.Range("I3").Value = "Yes"
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = $null
$ws.Range("O3").Value = $null

# Row 4 is removed entirely
$ws.Rows("4").Delete()

# Update selection to reflect the new active cell
$ws.Range("C4").Select()
